# PT Setting Fix
# - Bump the PT values on "PositiveExtra" from 10.5 to 12.5 (B2:B16)
# - Restore/park the selection on each sheet the way the author left them
# - Make "ProductList" the active (visible) tab when the workbook is reopened

$wb = $excel.ActiveWorkbook

# --- PositiveExtra ------------------------------------------------------
$wsPositiveExtra = $wb.Worksheets.Item("PositiveExtra")
$wsPositiveExtra.Activate()

# PT column values: 10.5 -> 12.5
$wsPositiveExtra.Range("B2:B16").Value = 12.5

# Leave the selection parked on the data that was just touched
$wsPositiveExtra.Range("B2:B16").Select()

# --- ProductList (left as the active/visible sheet) ---------------------
$wsProductList = $wb.Worksheets.Item("ProductList")
$wsProductList.Activate()
$wsProductList.Range("A16:B16").Select()

Write-Output "PT setting fix applied"
